# API Integration & Configuration test cases are now done: mark their
# "Passed" column as YES and stamp the completion date (1/23/2014) on the
# "Date" column, copying the existing date cell's number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$doneDate = (Get-Date -Year 2014 -Month 1 -Day 23 -Hour 0 -Minute 0 -Second 0).Date

# xlPasteFormats - used to stamp the existing "Date" column's date format
# (numFmtId 14) onto cells that were previously blank, without minting a
# brand-new number format.
$xlPasteFormats = -4122

# Row 3  -> ID 1  (APITest_Integration)
$ws.Range("D3").Value = "YES"
$ws.Range("E3").Value = $doneDate

# Row 14 -> ID 12 (APITest_Integration)
$ws.Range("D14").Value = "YES"
$ws.Range("E4").Copy()
$ws.Range("E14").PasteSpecial($xlPasteFormats)
$ws.Range("E14").Value = $doneDate

# Row 16 -> ID 14 (APITest_Integration)
$ws.Range("D16").Value = "YES"
$ws.Range("E4").Copy()
$ws.Range("E16").PasteSpecial($xlPasteFormats)
$ws.Range("E16").Value = $doneDate

$excel.CutCopyMode = $false

# Last cell the author left selected before saving.
$ws.Range("E16").Select()
